$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenarioA")

$ws.Range("F8").Formula = "=ROUND(0.6/G8,0)"
$ws.Range("G8").Value = 0.000405
$ws.Range("G8").NumberFormat = "0.000000"
$ws.Range("H8").Value = 0

$ws.Range("K9").Select()
